$d = $word.ActiveDocument

# Target the first paragraph (the **ID__...__ID** marker paragraph).
$para1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, 5pt space from text) and
# bump the left indent from 6pt (120 twips) to 11.25pt (225 twips).
$borders = $para1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
$para1.Range.ParagraphFormat.LeftIndent = 11.25

# Replace the placeholder id text (plus the trailing space run) with the
# updated id, collapsing the two runs into one. MatchWildcards is left
# off so the literal "**" in the search text is not treated as a glob.
$d.Content.Find.Execute("**ID__AFFARS_5313_topic_4__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5313_106_3__ID**", 2)
